$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 221, shifting existing rows 221:278 down to 222:279
$ws.Rows.Item(221).Insert()

# Populate the newly inserted row 221 with the new weekly record
$ws.Range("A221").Value = 8
$ws.Range("B221").Value = "Terminal La Palmera de La Serena"
$ws.Range("C221").Value = "Coquimbo"
$ws.Range("D221").Value = 44722
$ws.Range("E221").Value = 4
$ws.Range("F221").Value = 100112012
$ws.Range("G221").Value = "Espinaca"
$ws.Range("H221").Value = "Sin especificar"
$ws.Range("I221").Value = "Primera"
$ws.Range("J221").Value = 3000
$ws.Range("K221").Value = 500
$ws.Range("L221").Value = 600
$ws.Range("M221").Value = 550
$ws.Range("N221").Value = "$/atado 300 a 500 gramos"
$ws.Range("O221").Value = "Provincia del Elquí"
$ws.Range("P221").Value = 1100
$ws.Range("Q221").Value = 0.5
$ws.Range("R221").Value = "Hortaliza"
